$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Extend the "XXX Relatedly, ..." paragraph: add the sentence about
#    the disagreement with Mannion over the Xenoposeidon revision
#    (with the genus name italicised).
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("are a stronger focus for the introduction.")
$pt = $r.Duplicate
$pt.Collapse(0)
$pt.InsertAfter(" However, the disagreement with Mannion in the reviews of the ")
$pt.Collapse(0)

$pt.InsertAfter("Xenoposeidon")
$pt.Font.Italic = $true
$pt.Collapse(0)

$pt.InsertAfter(" revision remain of interest as a concrete example of how the ideas in the present manuscript are significant, so the account is retained in a reduced form in the discussion.")
$pt.Collapse(0)

# ------------------------------------------------------------------
# 2. Insert a brand new paragraph (about the Saegusa and Ikeda figure)
#    right after the paragraph just extended.
# ------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("so the account is retained in a reduced form in the discussion.")
$hostPara = $r2.Paragraphs(1)
$hostPara.Range.InsertParagraphAfter() | Out-Null

$r2b = $d.Content
$r2b.Find.Execute("so the account is retained in a reduced form in the discussion.")
$hostPara2b = $r2b.Paragraphs(1)
$newPara1 = $hostPara2b.Next()
$newPara1.Range.Text = "XXX We have regretfully removed the Saegusa and Ikeda’s (2014: figure 8) illustration, even though it provides a perfect single-figure example of the inconsistency we seek to illustrate, for fear of seeming confrontational. We have substituted less helpful but more inoffensive references to illustrations in our own work."

# ------------------------------------------------------------------
# 3. Extend the "XXX We have reworked and combined ..." paragraph with
#    the sentence about Figure 4 / Giraffatitan C5.
# ------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("along the lines that Marek suggests.")
$pt3 = $r3.Duplicate
$pt3.Collapse(0)
$pt3.InsertAfter(" We have also followed his suggestion to expand Figure 4 to show additional long cervical vertebrae as well as the ")
$pt3.Collapse(0)

$pt3.InsertAfter("Giraffatitan")
$pt3.Font.Italic = $true
$pt3.Collapse(0)

$pt3.InsertAfter(" C5 that is our principal example.")
$pt3.Collapse(0)

# ------------------------------------------------------------------
# 4. Insert three brand new paragraphs after the paragraph just
#    extended: colour-blindness, "In general, ...", and the open
#    peer-review paragraph (in that order).
# ------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("C5 that is our principal example.")
$hostPara3 = $r4.Paragraphs(1)
$hostPara3.Range.InsertParagraphAfter() | Out-Null

$r4b = $d.Content
$r4b.Find.Execute("C5 that is our principal example.")
$hostPara3b = $r4b.Paragraphs(1)
$newPara2 = $hostPara3b.Next()
$newPara2.Range.Text = "XXX We note the comments on colour-blindness, and have checked that where colour is used to distinguish lines with different meanings, those lines are also distinguishable by other criteria."

$r5 = $d.Content
$r5.Find.Execute("those lines are also distinguishable by other criteria.")
$hostPara4 = $r5.Paragraphs(1)
$hostPara4.Range.InsertParagraphAfter() | Out-Null

$r5b = $d.Content
$r5b.Find.Execute("those lines are also distinguishable by other criteria.")
$hostPara4b = $r5b.Paragraphs(1)
$newPara3 = $hostPara4b.Next()
$newPara3.Range.Text = "In general, we have followed the specific comments attached to Marek’s review, but with some exceptions. In particular, we do not agree with the suggestion that expunging the pronoun “we” throughout and substituting passive voice would improve the manuscript. We have changed some instances, but in other clarity is best served by the simpler active phrasing."

$r6 = $d.Content
$r6.Find.Execute("is best served by the simpler active phrasing.")
$hostPara5 = $r6.Paragraphs(1)
$hostPara5.Range.InsertParagraphAfter() | Out-Null

$r6b = $d.Content
$r6b.Find.Execute("is best served by the simpler active phrasing.")
$hostPara5b = $r6b.Paragraphs(1)
$newPara4 = $hostPara5b.Next()
$newPara4.Range.Text = "We have retained the section on open peer-review, contrary to Marek’s recommendation, as it is important to us and relevant to the origin of the present paper’s core question."
